# Updated cryptos list — price (D) and 1h volume-change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.200.40'
$ws.Range("E2").Value = '  +1.30%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.649.29'
$ws.Range("E3").Value = '  -0.80%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '607.91'
$ws.Range("E5").Value = '  -0.55%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.69'
$ws.Range("E6").Value = '  +3.65%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("E8").Value = '  +0.46%  '

$ws.Range("E9").Value = '  +1.89%  '

$ws.Range("E10").Value = '  +6.43%  '

$ws.Range("E11").Value = '  -0.18%  '

$ws.Range("E12").Value = '  -0.91%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.56'
$ws.Range("E13").Value = '  +0.88%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.124.13'
$ws.Range("E14").Value = '  -0.84%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.036.68'
$ws.Range("E15").Value = '  +1.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000148'
$ws.Range("E16").Value = '  +1.99%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.635.27'
$ws.Range("E17").Value = '  -1.65%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.92'
$ws.Range("E18").Value = '  +4.35%  '

$ws.Range("E19").Value = '  +3.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '347.02'
$ws.Range("E20").Value = '  +1.55%  '

$ws.Range("E21").Value = '  +0.19%  '

$ws.Range("E22").Value = '  -0.03%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.58'
$ws.Range("E23").Value = '  -0.44%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.30'
$ws.Range("E24").Value = '  -1.20%  '

$ws.Range("E25").Value = '  +8.05%  '

$ws.Range("E26").Value = '  +4.84%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.16'
$ws.Range("E28").Value = '  +3.44%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '552.81'
$ws.Range("E29").Value = '  +1.70%  '

$ws.Range("E30").Value = '  -1.81%  '

$ws.Range("E31").Value = '  -0.12%  '

$ws.Range("E32").Value = '  +0.00%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0852'
$ws.Range("E33").Value = '  +5.88%  '

$ws.Range("E34").Value = '  -2.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.27'
$ws.Range("E35").Value = '  +2.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '169.49'
$ws.Range("E36").Value = '  -2.05%  '

$ws.Range("E37").Value = '  +0.54%  '

$ws.Range("E38").Value = '  +0.01%  '

$ws.Range("E39").Value = '  +4.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '19.31'
$ws.Range("E40").Value = '  +0.70%  '

$ws.Range("E41").Value = '  +0.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '165.85'
$ws.Range("E42").Value = '  -6.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '40.40'
$ws.Range("E43").Value = '  +0.73%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.82'
$ws.Range("E44").Value = '  +2.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.10'
$ws.Range("E45").Value = '  -0.62%  '

$ws.Range("E46").Value = '  -1.07%  '

$ws.Range("E47").Value = '  -0.97%  '

$ws.Range("E48").Value = '  +14.33%  '

$ws.Range("E49").Value = '  +2.30%  '

$ws.Range("E50").Value = '  -0.20%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '18.91'
$ws.Range("E51").Value = '  +1.06%  '
